# Update BOM with Ebay item numbers for critical parts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "JST XH series knockoff" references with a footnote asterisk
$ws.Range("I7").Value  = "JST XH series knockoff*"
$ws.Range("I9").Value  = "JST XH series knockoff*. Remove center pin."
$ws.Range("I26").Value = "JST XH series knockoff*"
$ws.Range("I27").Value = "JST XH series knockoff*"
$ws.Range("I28").Value = "JST XH series knockoff*"

# Add Ebay item numbers for critical parts (J6/J7 connector, K1 relay, PSU)
$ws.Range("H10").Value = "Ebay 231260624408"
$ws.Range("H11").Value = "Ebay 361282277725"
$ws.Range("H25").Value = "Ebay 351211320656"

# Tidy up note punctuation
$ws.Range("I25").Value = "see picture, hole pattern must match main board."

# Add footnote explaining the asterisk, referencing the Ebay kit listing
$ws.Range("B35").Value = "*"
$ws.Range("D35").Value = "A small kit of JST XH compatible"
$ws.Range("D36").Value = "connectors is available on Ebay"
$ws.Range("D37").Value = "under item 141292096528"
$ws.Range("D38").Value = "This will have all of the headers"
$ws.Range("D39").Value = "pins and housings you need"
$ws.Range("D40").Value = "to build the project."

# Widen Source column slightly to fit the new Ebay item numbers
$ws.Columns.Item(8).ColumnWidth = 18.1

# Move the active selection to reflect where editing left off
[void]$ws.Range("D41").Select()
